$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C ("Förändrad" date) for rows 2 through 45 from serial 45233 to 45243
for ($row = 2; $row -le 45; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45233) {
        $cell.Value2 = 45243
    }
}
